$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.740127682685852
$ws.Range("B1").Value = 2.297441720962524
$ws.Range("C1").Value = 4.628640651702881
$ws.Range("D1").Value = 4.179306983947754
$ws.Range("E1").Value = 1.682458162307739
